$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing values
$ws.Range("B2").Value = 5
$ws.Range("B3").Value = 4.0999999999999996
$ws.Range("B4").Value = 0.8

# Row 5 becomes the new "theta_threshold_range" row
$ws.Range("A5").Value = "theta_threshold_range"
$ws.Range("B5").Value = 1.5
$ws.Range("C5").Value = 140

# Row 6 gets the former "pie_threshold_range" row (moved down)
$ws.Range("A5:C5").Copy()
$ws.Range("A6:C6").PasteSpecial(-4122)

$ws.Range("A6").Value = "pie_threshold_range"
$ws.Range("B6").Value = 0
$ws.Range("C6").Value = 19.899999999999999

$ws.Range("B6").Select()
